$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "2025-09-19 Friday" "2025-09-20 Saturday"

Replace-Text "170×9=1530" "936×7=6552"
Replace-Text "827×7=5789" "184×8=1472"
Replace-Text "387×7=2709" "864×3=2592"
Replace-Text "106×3=318" "808×9=7272"
Replace-Text "753×4=3012" "681×5=3405"

Replace-Text "326×9=2934" "382×2=764"
Replace-Text "292×4=1168" "181×2=362"
Replace-Text "779×4=3116" "337×6=2022"
Replace-Text "215×5=1075" "401×7=2807"
Replace-Text "425×9=3825" "549×2=1098"

Replace-Text "994×5=4970" "657×7=4599"
Replace-Text "950×4=3800" "760×7=5320"
Replace-Text "424×2=848" "446×9=4014"
Replace-Text "472×2=944" "933×6=5598"
Replace-Text "168×3=504" "357×7=2499"

Replace-Text "676×7=4732" "931×2=1862"
Replace-Text "932×6=5592" "385×6=2310"
Replace-Text "103×5=515" "766×5=3830"
Replace-Text "206×9=1854" "477×7=3339"
Replace-Text "353×6=2118" "287×4=1148"

Replace-Text "261×9=2349" "943×4=3772"
Replace-Text "702×2=1404" "974×7=6818"
Replace-Text "983×4=3932" "687×9=6183"
Replace-Text "333×9=2997" "446×3=1338"
Replace-Text "244×4=976" "830×6=4980"
